$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2938.6667

$ws.Range("H26").Value = 15347.667
$ws.Range("I26").Value = 10013
$ws.Range("J26").Value = 18015
$ws.Range("K26").Value = 10013
$ws.Range("L26").Value = 18015
$ws.Range("M26").Value = -9669
$ws.Range("N26").Value = -18703

$ws.Range("H132").Value = 5937.778
$ws.Range("I132").Value = 6070.625
$ws.Range("J132").Value = 4875
$ws.Range("K132").Value = 18211.875
$ws.Range("L132").Value = 14625
$ws.Range("M132").Value = -15681.875
$ws.Range("N132").Value = -19685

$ws.Range("H138").Value = 2171.0435
$ws.Range("J138").Value = 2999.9333
$ws.Range("L138").Value = 8999.7999
$ws.Range("N138").Value = -19279.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = -856
$ws.Range("N13").Value = -4288

$ws.Range("H14").Value = 1200
$ws.Range("I14").Value = 1200
$ws.Range("K14").Value = 1200
$ws.Range("M14").Value = -1025

$ws.Range("H32").Value = 3790.8333
$ws.Range("I32").Value = 3992.3333
$ws.Range("K32").Value = 3992.3333
$ws.Range("M32").Value = -3705.3333

$ws.Range("H45").Value = 900
$ws.Range("I45").Value = 900
$ws.Range("K45").Value = 900
$ws.Range("M45").Value = -523

$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 20000
$ws.Range("K55").Value = 20000
$ws.Range("M55").Value = -19685

$ws.Range("H63").Value = 5891.6665
$ws.Range("I63").Value = 4070
$ws.Range("J63").Value = 15000
$ws.Range("K63").Value = 4070
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -3384
$ws.Range("N63").Value = -16372

$ws.Range("H66").Value = 5891.6665
$ws.Range("I66").Value = 4070
$ws.Range("J66").Value = 15000
$ws.Range("K66").Value = 20350
$ws.Range("L66").Value = 75000
$ws.Range("M66").Value = -16918
$ws.Range("N66").Value = -81864

$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

$ws.Range("H132").Value = 1903
$ws.Range("I132").Value = 1204
$ws.Range("K132").Value = 3612
$ws.Range("M132").Value = -1082

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13614.286
$ws.Range("J26").Value = 12450
$ws.Range("L26").Value = 12450
$ws.Range("N26").Value = -13034

$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 50000
$ws.Range("L35").Value = 50000
$ws.Range("N35").Value = -50620

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38
$ws.Range("I7").Value = 19.2
$ws.Range("J7").Value = 69.333336
$ws.Range("K7").Value = 19.2
$ws.Range("L7").Value = 69.333336
$ws.Range("M7").Value = 93.8
$ws.Range("N7").Value = -295.333336

$ws.Range("H16").Value = 499
$ws.Range("I16").Value = 499
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 499
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -212

$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H58").Value = 3884.8096
$ws.Range("I58").Value = 1916.4
$ws.Range("K58").Value = 1916.4
$ws.Range("M58").Value = -1713.4

$ws.Range("H113").Value = 499
$ws.Range("I113").Value = 499
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 499
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1671

$ws.Range("H116").Value = 35000
$ws.Range("J116").Value = 35000
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H132").Value = 6473.647
$ws.Range("I132").Value = 2831.5
$ws.Range("K132").Value = 8494.5
$ws.Range("M132").Value = -5964.5

$ws.Range("H134").Value = 2566.4119
$ws.Range("I134").Value = 1362.9
$ws.Range("K134").Value = 4088.7
$ws.Range("M134").Value = -1553.7

$ws.Range("H136").Value = 3884.8096
$ws.Range("I136").Value = 1916.4
$ws.Range("K136").Value = 5749.200000000001
$ws.Range("M136").Value = -3199.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H6").Value = 90952.45
$ws.Range("I6").Value = 41.666668
$ws.Range("J6").Value = 200045.4
$ws.Range("K6").Value = 125.000004
$ws.Range("L6").Value = 600136.2
$ws.Range("M6").Value = -12.000004
$ws.Range("N6").Value = -600362.2

$ws.Range("H21").Value = 10
$ws.Range("I21").Value = 10
$ws.Range("K21").Value = 30
$ws.Range("M21").Value = 143

$ws.Range("H23").Value = 373.5
$ws.Range("I23").Value = 291
$ws.Range("K23").Value = 873
$ws.Range("M23").Value = -638

$ws.Range("H26").Value = 64.85714
$ws.Range("I26").Value = 50.666668
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 152.000004
$ws.Range("L26").Value = 450
$ws.Range("M26").Value = 135.999996
$ws.Range("N26").Value = -1026

$ws.Range("H48").Value = 100
$ws.Range("I48").Value = 100
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 300
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -50

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 9095
$ws.Range("J54").Value = 9095
$ws.Range("L54").Value = 9095
$ws.Range("N54").Value = -9875

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492

$ws.Range("H97").Value = 1552.375
$ws.Range("I97").Value = 1552.375
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1552.375
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1056.375

$ws.Range("H102").Value = 1442.8572
$ws.Range("I102").Value = 1442.8572
$ws.Range("K102").Value = 1442.8572
$ws.Range("M102").Value = 179.1428000000001

$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2196.3333
$ws.Range("I16").Value = 2237.6
$ws.Range("K16").Value = 2237.6
$ws.Range("M16").Value = -2067.6

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws.Range("H55").Value = 1234
$ws.Range("I55").Value = 1234
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1234
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1061

$ws.Range("H122").Value = 5638
$ws.Range("I122").Value = 5797.5
$ws.Range("K122").Value = 17392.5
$ws.Range("M122").Value = -14942.5

$ws.Range("H136").Value = 11000
$ws.Range("I136").Value = 11000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 33000
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -30450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 200
$ws.Range("I13").Value = 200
$ws.Range("K13").Value = 200
$ws.Range("M13").Value = -60

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 1071.2858
$ws.Range("I107").Value = 350
$ws.Range("K107").Value = 1050
$ws.Range("M107").Value = 870

$ws.Range("H132").Value = 5819.1
$ws.Range("I132").Value = 6024
$ws.Range("K132").Value = 18072
$ws.Range("M132").Value = -15542

$ws.Range("H136").Value = 1769.6
$ws.Range("I136").Value = 1769.6
$ws.Range("K136").Value = 5308.799999999999
$ws.Range("M136").Value = -2758.799999999999
